$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from C1 into the two new header cells D1:E1
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null

# Header row text
$ws.Range("A1").Value = "Age Category"
$ws.Range("B1").Value = "Purchase Count"
$ws.Range("C1").Value = "Average Purchase Price"
$ws.Range("D1").Value = "Total Purchase Value"
$ws.Range("E1").Value = "Average Total Purchase per Person"

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2 (<10)
$ws.Range("B2").Value = 23
Set-TextCell $ws.Range("C2") "$3.35"
Set-TextCell $ws.Range("D2") "$77.13"
Set-TextCell $ws.Range("E2") "$0.13"

# Row 3 (10-14)
$ws.Range("B3").Value = 28
Set-TextCell $ws.Range("C3") "$2.96"
Set-TextCell $ws.Range("D3") "$82.78"
Set-TextCell $ws.Range("E3") "$0.14"

# Row 4 (15-19)
$ws.Range("B4").Value = 136
Set-TextCell $ws.Range("C4") "$3.04"
Set-TextCell $ws.Range("D4") "$412.89"
Set-TextCell $ws.Range("E4") "$0.72"

# Row 5 (20-24)
$ws.Range("B5").Value = 365
Set-TextCell $ws.Range("C5") "$3.05"
Set-TextCell $ws.Range("D5") "$1,114.06"
Set-TextCell $ws.Range("E5") "$1.93"

# Row 6 (25-29)
$ws.Range("B6").Value = 101
Set-TextCell $ws.Range("C6") "$2.90"
Set-TextCell $ws.Range("D6") "$293.00"
Set-TextCell $ws.Range("E6") "$0.51"

# Row 7 (30-34)
$ws.Range("B7").Value = 73
Set-TextCell $ws.Range("C7") "$2.93"
Set-TextCell $ws.Range("D7") "$214.00"
Set-TextCell $ws.Range("E7") "$0.37"

# Row 8 (35-39)
$ws.Range("B8").Value = 41
Set-TextCell $ws.Range("C8") "$3.60"
Set-TextCell $ws.Range("D8") "$147.67"
Set-TextCell $ws.Range("E8") "$0.26"

# Row 9 (40+)
$ws.Range("B9").Value = 13
Set-TextCell $ws.Range("C9") "$2.94"
Set-TextCell $ws.Range("D9") "$38.24"
Set-TextCell $ws.Range("E9") "$0.07"

$wb.Save()
